$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: card holder name / card number ---
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 13.12.2024"

# --- Row 6 (transaction 1) ---
$ws.Range("B6").Value = "14.12."
$ws.Range("C6").Value = "15.12."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 94045574"
$ws.Range("E6").Value = "84,95-"

# --- Row 7 (transaction 2) ---
$ws.Range("B7").Value = "18.12."
$ws.Range("C7").Value = "19.12."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "24,82-"

# --- Row 8 (transaction 3) ---
$ws.Range("B8").Value = "22.12."
$ws.Range("C8").Value = "23.12."
$ws.Range("D8").Value = "PAYPAL WEHURX"
$ws.Range("E8").Value = "22,14-"

# --- Row 9: no longer used for a transaction - clear it out and adopt the
#     "blank row" formatting (E column centered, like the other blank rows) ---
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# --- Row 10: no longer used for a transaction - clear it out and adopt the
#     "blank row" formatting (E column right aligned, like row 11) ---
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 27.12.2024"
$ws.Range("E12").Value = "131,91-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 05.01.2025"
